$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Force text format on cells whose new values look numeric, to preserve original string formatting
$numericRefs = @('D5', 'D6', 'D9', 'D10', 'D12', 'D16', 'D17', 'D19', 'D20', 'D21', 'D22', 'D23', 'D26', 'D27', 'D28', 'D29', 'D30', 'D31', 'D33', 'D36', 'D37', 'D38', 'D40', 'D41', 'D43', 'D44', 'D45', 'D46', 'D47', 'D49', 'D50', 'D51')
foreach ($r in $numericRefs) {
    $ws.Range($r).NumberFormat = "@"
}

$ws.Range('D2').Value = '60.404.80'
$ws.Range('E2').Value = '  -4.41%  '
$ws.Range('D3').Value = '3.348.96'
$ws.Range('E3').Value = '  -2.28%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '566.95'
$ws.Range('E5').Value = '  -2.38%  '
$ws.Range('D6').Value = '132.21'
$ws.Range('E6').Value = '  +2.43%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = '3.344.04'
$ws.Range('E8').Value = '  -2.40%  '
$ws.Range('D9').Value = '0.473'
$ws.Range('E9').Value = '  -1.18%  '
$ws.Range('D10').Value = '7.60'
$ws.Range('E10').Value = '  +0.70%  '
$ws.Range('E11').Value = '  -3.26%  '
$ws.Range('D12').Value = '0.379'
$ws.Range('E12').Value = '  -0.44%  '
$ws.Range('D13').Value = '3.918.85'
$ws.Range('E13').Value = '  -2.37%  '
$ws.Range('E14').Value = '  +0.12%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '3.348.90'
$ws.Range('E15').Value = '  -2.35%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').Value = '0.0000171'
$ws.Range('E16').Value = '  -2.82%  '
$ws.Range('D17').Value = '24.69'
$ws.Range('E17').Value = '  -1.49%  '
$ws.Range('D18').Value = '60.466.56'
$ws.Range('E18').Value = '  -4.34%  '
$ws.Range('D19').Value = '13.61'
$ws.Range('E19').Value = '  +2.74%  '
$ws.Range('B20').Value = 'Polkadot'
$ws.Range('C20').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D20').Value = '5.74'
$ws.Range('E20').Value = '  +0.90%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').Value = '9.28'
$ws.Range('E21').Value = '  -5.21%  '
$ws.Range('D22').Value = '373.84'
$ws.Range('E22').Value = '  -2.40%  '
$ws.Range('D23').Value = '0.560'
$ws.Range('E23').Value = '  -0.86%  '
$ws.Range('D24').Value = '3.480.21'
$ws.Range('E24').Value = '  -2.38%  '
$ws.Range('E25').Value = '  +0.19%  '
$ws.Range('D26').Value = '69.52'
$ws.Range('E26').Value = '  -5.21%  '
$ws.Range('D27').Value = '0.0000114'
$ws.Range('E27').Value = '  +3.81%  '
$ws.Range('D28').Value = '1.63'
$ws.Range('E28').Value = '  +16.06%  '
$ws.Range('D29').Value = '7.50'
$ws.Range('E29').Value = '  +5.99%  '
$ws.Range('D30').Value = '0.999'
$ws.Range('E30').Value = '  -0.03%  '
$ws.Range('D31').Value = '8.04'
$ws.Range('E31').Value = '  +1.55%  '
$ws.Range('E32').Value = '  +0.04%  '
$ws.Range('D33').Value = '2.13'
$ws.Range('E33').Value = '  -2.70%  '
$ws.Range('E34').Value = '  -0.05%  '
$ws.Range('D35').Value = '3.379.64'
$ws.Range('E35').Value = '  -2.21%  '
$ws.Range('D36').Value = '22.94'
$ws.Range('E36').Value = '  +1.06%  '
$ws.Range('D37').Value = '5.40'
$ws.Range('E37').Value = '  +1.78%  '
$ws.Range('D38').Value = '6.92'
$ws.Range('E38').Value = '  +2.24%  '
$ws.Range('E39').Value = '  -0.16%  '
$ws.Range('D40').Value = '159.94'
$ws.Range('E40').Value = '  -2.62%  '
$ws.Range('D41').Value = '0.0773'
$ws.Range('E41').Value = '  +0.71%  '
$ws.Range('E42').Value = '  -0.14%  '
$ws.Range('D43').Value = '1.20'
$ws.Range('E43').Value = '  +9.34%  '
$ws.Range('D44').Value = '4.38'
$ws.Range('E44').Value = '  +1.25%  '
$ws.Range('D45').Value = '41.10'
$ws.Range('E45').Value = '  -0.07%  '
$ws.Range('D46').Value = '0.752'
$ws.Range('E46').Value = '  -3.77%  '
$ws.Range('D47').Value = '23.88'
$ws.Range('E47').Value = '  +2.71%  '
$ws.Range('E48').Value = '  -0.87%  '
$ws.Range('D49').Value = '6.87'
$ws.Range('E49').Value = '  +2.08%  '
$ws.Range('D50').Value = '22.62'
$ws.Range('E50').Value = '  +10.96%  '
$ws.Range('D51').Value = '0.893'
$ws.Range('E51').Value = '  +0.73%  '

# Restore default style on the numeric-look cells so no leftover style index remains
foreach ($r in $numericRefs) {
    $ws.Range($r).Style = "Normal"
}
